$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: "Rotate Array to right by K" (GFG) gets a LeetCode number (189)
# and both cells get a yellow highlight fill ---
$ws.Range("A17").Value = 189
$ws.Range("A17").Style = "Normal"
$ws.Range("A17").Interior.Color = 65535
$ws.Range("A17").HorizontalAlignment = -4108

$ws.Range("B17").Style = "Normal"
$ws.Range("B17").Interior.Color = 65535

# --- New row 23: "Rotate Array  to left by K" (GFG / Java), same highlight ---
$ws.Range("A23").Value = "GFG"
$ws.Range("A23").Style = "Normal"
$ws.Range("A23").Interior.Color = 65535
$ws.Range("A23").HorizontalAlignment = -4108

$ws.Range("B23").Value = "Rotate Array  to left by K"
$ws.Range("B23").Style = "Normal"
$ws.Range("B23").Interior.Color = 65535

$ws.Range("C23").Value = "Java"

$ws.Range("D23").Value = 45000
$ws.Range("D23").NumberFormat = $ws.Range("D22").NumberFormat

# --- View state: scroll so row 4 is at top, select J20 ---
$app = $excel
$app.Goto($ws.Range("A4"), $true)
$ws.Range("J20").Select()
